# TB_BOM.xlsx update: add a "Front Shoulder" measurement column to the
# parts/lengths table, per commit "updated BOM with shoulder".
#
# The new column is inserted right after "Front Leg" (column H) and before
# the old "Front Foot" column, pushing every later column (old I..R) one
# slot to the right (new J..S). Row 2 gets a shoulder count of 1, row 4
# gets a shoulder count of 7; the totals/sum formulas that summed the
# H:K(old)/H:L(new) leg+foot block and the grand totals automatically pick
# up the new column because they are re-entered against the shifted ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at I (old I "Front Foot" and everything to its
# right - Back Leg, Back Foot, x2, Sum, Add 10%, Total Boards, Total board +
# 10% - shift right by one column: I->J, J->K, K->L, L->M, N->O, O->P,
# Q->R, R->S). Formulas that reference these cells are shifted/expanded
# along with them.
$ws.Columns("I").Insert() | Out-Null

# New column header + data for "Front Shoulder".
$ws.Range("I1").Value = "Front Shoulder"
$ws.Range("I2").Value = 1
$ws.Range("I4").Value = 7

# Re-assert the "x2" subtotal column (now column M, was L) so it includes
# the new Front Shoulder column in its sum, and so the un-touched rows
# (3-7) stay one shared formula like the original layout.
$ws.Range("M2").Formula = "=SUM(H2:L2)*2"
$ws.Range("M3:M7").Formula = "=SUM(H3:L3)*2"

# Re-assert the "Add 10%" column (now column P, was O) as a single shared
# formula over its un-touched rows, matching the original layout shape.
$ws.Range("P3:P6").Formula = "=SUM(O3+O3*0.1)"

# Leave the selection where the author left it when they saved.
$ws.Range("J12").Select() | Out-Null
